$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3 ---

# Tanggal_Jatuh_Tempo: numeric date serial, formatted as YYYY-MM-DD.
# Done first so the new custom-format style claims cellXfs index 2,
# matching the target workbook's style table ordering.
$ws.Range("F3").Value = 46240
$ws.Range("F3").NumberFormat = "yyyy-mm-dd"
$ws.Range("F3").NumberFormat = "YYYY-MM-DD"

# NIK (A3) looks like a pure-digit string but must stay text, like A2.
$ws.Range("A3").Value = "'1234456278949533"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "BG6744HU"
$ws.Range("C3").Value = "Rahma"
$ws.Range("D3").Value = "Palembang"

# Pajak_Terhutang / Pajak are numeric this time (unlike row 2's text values)
$ws.Range("E3").Value = 65000
$ws.Range("G3").Value = 65000

$ws.Range("H3").Value = "HIDJR3544H"
$ws.Range("I3").Value = "Yamaha"
$ws.Range("J3").Value = "Sepeda Motor"
$ws.Range("K3").Value = "Hitam"

# Status / Status_Pengiriman / No_Resi / Ekspedisi left blank for this entry
$ws.Range("L3").Value = "'"
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Value = "'"
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Value = "'"
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value = "'"
$ws.Range("O3").Style = "Normal"
